$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "JSU(-1.2840078650272977, 1.4202899036264873, -0.15642304795461545, 3.5236745814291814)"
$ws.Range("C2").Value = "JSU(-1.2791876249815466, 1.0161284618576734, 2.582774038851654, 3.6756187909441014)"
$ws.Range("D2").Value = "JSU(-0.8408061581559201, 1.1123423849647511, 0.6150053495213439, 2.6819460436439813)"
$ws.Range("E2").Value = "NIG(0.8581562897097794, 0.5868966586878821, 5.443565967109327, 5.349445340741092)"
